$d = $word.ActiveDocument

function Replace-Text($findText, $replaceText) {
    $d.Content.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)
}

# "User Acceptance Tesing" -> "User Acceptance Testing"
Replace-Text "Tesing" "Testing"

# Tidy the acceptance-criteria paragraph: supervisors->supervisor, double spaces, Wach->Each
Replace-Text "my project supervisors for the functionality that  I plan on implementing for that week.  These" "my project supervisor for the functionality that I plan on implementing for that week. These"
Replace-Text "the user stories outlined in the project plan.  Wach set of unit tests" "the user stories outlined in the project plan. Each set of unit tests"

# Unit testing paragraph tidy-ups
Replace-Text "a good indication  as to whether" "a good indication as to whether"
Replace-Text "I will be using the Scalatest unit testing tool to do this which runs on top of JUnit ( the standard Java tool for unit testing)." "I will be using the ScalaTest unit testing tool to do this, which runs on top of JUnit (the standard Java tool for unit testing)."
Replace-Text "how I document this section of testing.  I will also be using the standard scala test coverage tool" "how I document this section of testing. I will also be using the standard Scala test coverage tool"

# Integration testing paragraph tidy-ups
Replace-Text "work as they should. I will preform extensive manual testing of the system in order to ensure that the new functionalities work as they should. These tests will be more in deph versions of the user acceptance tests that I specified at the start of each iteration , with the aim being that  by the end of this phase of testing,  should be confident that the functionality  that I have implemented works and can be presented to my project supervisor t the next progress meeting.  This phase" "work as they should, I will perform extensive manual testing of the system in order to ensure that the new functionalities work as they should. These tests will be more in depth versions of the user acceptance tests that I specified at the start of each iteration, with the aim being that by the end of this phase of testing, I should be confident that the functionality that I have implemented works and can be presented to my project supervisor at the next progress meeting. This phase"

# Make the "Integration Testing" heading bold, matching the other section headings
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.Trim() -eq "Integration Testing") {
        $p.Range.Bold = 1
    }
}
